$d = $word.ActiveDocument

# 1. Fix typo "cancell" -> "cancel" in the Manager-cancels-option sentence.
$d.Content.Find.Execute(
    "The Manager chooses the cancell option:", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "The Manager chooses the cancel option:", 2) | Out-Null

# 2. Word auto-tracks the most recent edit location with the single
#    "_GoBack" bookmark. Re-add it at the point right after "cancel"
#    (collapsed, zero-length) -- this both relocates it from its old spot
#    (after ": Sometimes") and removes the stale one, since a document can
#    only have one bookmark with a given name.
$r = $d.Content
$r.Find.Execute("cancel", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

# 3. Update the cached SAVEDATE field text in the footer.
$ftr = $d.Sections(1).Footers(1)
$ftr.Range.Find.Execute(
    "2/21/2017 5:42:00 PM", $true, $false, $false, $false, $false, $true,
    1, $false, "2/21/2017 8:56:00 PM", 2) | Out-Null

# 4. Update the cached DATE field text in the header.
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Find.Execute(
    "2/21/2017", $true, $false, $false, $false, $false, $true, 1, $false,
    "2/22/2017", 2) | Out-Null
